$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.191.30"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").Value = "3.559.36"
$ws.Range("E3").Value = "  +1.99%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'605.47"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D6").Value = "'144.42"
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("D7").Value = "3.558.47"
$ws.Range("E7").Value = "  +2.06%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  +1.06%  "
$ws.Range("E10").Value = "  -3.43%  "
$ws.Range("D11").Value = "'8.04"
$ws.Range("E11").Value = "  +1.77%  "
$ws.Range("E12").Value = "  -1.92%  "
$ws.Range("D13").Value = "4.165.93"
$ws.Range("E13").Value = "  +2.18%  "
$ws.Range("E14").Value = "  -2.37%  "
$ws.Range("D15").Value = "'30.25"
$ws.Range("E15").Value = "  -2.76%  "
$ws.Range("D16").Value = "3.557.18"
$ws.Range("E16").Value = "  +1.79%  "
$ws.Range("D17").Value = "66.291.08"
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("D18").Value = "'11.65"
$ws.Range("E18").Value = "  +8.80%  "
$ws.Range("E19").Value = "  -1.07%  "
$ws.Range("E20").Value = "  -1.28%  "
$ws.Range("D21").Value = "'14.96"
$ws.Range("E21").Value = "  -2.33%  "
$ws.Range("D22").Value = "'428.34"
$ws.Range("E22").Value = "  -0.65%  "
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("D24").Value = "'78.70"
$ws.Range("E24").Value = "  -0.94%  "
$ws.Range("D25").Value = "3.700.30"
$ws.Range("E25").Value = "  +2.10%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").Value = "'0.0000120"
$ws.Range("E27").Value = "  +3.39%  "
$ws.Range("D28").Value = "'8.07"
$ws.Range("E28").Value = "  -0.68%  "
$ws.Range("D29").Value = "'9.24"
$ws.Range("E29").Value = "  -5.13%  "
$ws.Range("D30").Value = "'2.49"
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("D32").Value = "'1.49"
$ws.Range("E32").Value = "  -3.95%  "
$ws.Range("E33").Value = "  -3.86%  "
$ws.Range("B34").Value = "RenzoRestakedETH"
$ws.Range("C34").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D34").Value = "3.556.53"
$ws.Range("E34").Value = "  +2.24%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").Value = "'25.44"
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  -1.20%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'5.66"
$ws.Range("E38").Value = "  -1.09%  "
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "'7.85"
$ws.Range("E39").Value = "  -1.05%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("D41").Value = "'171.79"
$ws.Range("E41").Value = "  -1.21%  "
$ws.Range("D42").Value = "'0.0859"
$ws.Range("E42").Value = "  -3.96%  "
$ws.Range("D43").Value = "'5.30"
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("D44").Value = "'0.894"
$ws.Range("E44").Value = "  +0.27%  "
$ws.Range("E45").Value = "  -4.04%  "
$ws.Range("D46").Value = "'45.86"
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").Value = "'1.21"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'26.00"
$ws.Range("E48").Value = "  -5.94%  "
$ws.Range("E49").Value = "  +0.41%  "
$ws.Range("E50").Value = "  -2.15%  "
$ws.Range("D51").Value = "'0.950"
$ws.Range("E51").Value = "  -3.20%  "
